$d = $word.ActiveDocument

# Update the date heading
[void]$d.Content.Find.Execute("2023-09-19 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-20 Wednesday", 2)

# Update the division exercise table. Each filled row/column is addressed
# directly by position so that values which coincide with other rows'
# old/new text (e.g. "16÷5=3, 1", "48÷4=12, 0") are never double-replaced.
$tbl = $d.Tables.Item(1)

$rowsData = @{
    1  = @("81÷2=40, 1", "34÷8=4, 2", "98÷3=32, 2", "30÷3=10, 0", "72÷2=36, 0")
    5  = @("12÷9=1, 3", "61÷4=15, 1", "23÷9=2, 5", "28÷2=14, 0", "95÷3=31, 2")
    9  = @("46÷4=11, 2", "16÷5=3, 1", "66÷2=33, 0", "60÷2=30, 0", "15÷5=3, 0")
    13 = @("69÷8=8, 5", "13÷7=1, 6", "37÷4=9, 1", "19÷3=6, 1", "31÷8=3, 7")
    17 = @("48÷4=12, 0", "20÷9=2, 2", "81÷8=10, 1", "27÷3=9, 0", "98÷9=10, 8")
}

foreach ($rowIndex in $rowsData.Keys) {
    $values = $rowsData[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
